$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the D1 header text (identity text changed)
$ws.Range("D1").Value = "身分 (學士、碩士或博士班）"

# Widen column D to fit the new, longer text.
# NOTE: Excel's COM ColumnWidth is quantized to whole pixels (MDW-based grid,
# increments of 1/7 character-width units here), so the raw OOXML <col width>
# of 29.125 isn't directly reachable as a ColumnWidth input - the nearest
# representable value snaps to 29.1428571428571..., which is what we target.
$ws.Columns("D").ColumnWidth = 28.4285714285714

# Move the active selection to D6
$ws.Range("D6").Select()
